$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Actualiza base de datos EC" - the account-statement rows (16-19) are
# re-sorted by "Periodo Mora" (ascending: 1604, 1605, 1606, 1608), keeping
# each period's "Valor Mora" paired exactly as before the sort.
$ws.Range("E16").Value = "1604"
$ws.Range("F16").Value = 27600

$ws.Range("E17").Value = "1605"
$ws.Range("F17").Value = 27600

$ws.Range("E18").Value = "1606"
$ws.Range("F18").Value = 27600

$ws.Range("E19").Value = "1608"
$ws.Range("F19").Value = 9193
